# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" (Overview!E2/F2,
#   zh-cn!C2, de-de!C2 all share the same underlying text).
# - The handoff timestamps are refreshed to the new generation time.
# - The now-wider "Ready for handoff" label needs a wider status column on
#   every sheet that shows it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps --------------------------------------------------------
# Overview!G2 and de-de!H2 shared the same "Latest HO Xliff Generate Date" /
# "Latest Handoff Datetime" value before the edit, and still match after.
$overview.Range("G2").Value = "2016-08-25 20:37:42"
$dede.Range("H2").Value     = "2016-08-25 20:37:42"

# zh-cn!H2 ("Latest Handoff Datetime") advances independently.
$zhcn.Range("H2").Value = "2016-08-25 20:37:37"

# --- Column widths -------------------------------------------------------
# Widen the status columns to fit "Ready for handoff". ColumnWidth is stored
# in whole-pixel steps, so 16.333333333333 is the closest achievable setting
# to the authored width.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333
$overview.Columns.Item(6).ColumnWidth = 16.333333333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.333333333333
$dede.Columns.Item(3).ColumnWidth     = 16.333333333333
